# Project Hephaestus design doc:
# Split the last paragraph ("Force connection: Craft piercing tools and
# jewellery.") so that two new paragraphs follow it: a blank line and a
# note for Jack. The "_GoBack" bookmark (Word's "last edit location"
# marker) must end up at the very end of the document, after the new
# note, instead of where it used to sit (right before "Force connection").

$d = $word.ActiveDocument

# 1) Grow the final paragraph into three paragraphs: the existing text,
#    a blank paragraph, then the new note -- tagging the very end of the
#    note with a temporary marker "~" so we can find it precisely even
#    while it is still the last character in the story.
$d.Content.Find.Execute(
    "Force connection: Craft piercing tools and jewellery.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Force connection: Craft piercing tools and jewellery.^p^pJack needs to find out how to create a UI~",
    2) | Out-Null

# 2) The "_GoBack" bookmark currently sits right before "Force
#    connection..."; drop it so we can recreate it in the right spot.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3) Re-create the bookmark collapsed right after "...UI", i.e. just
#    before the temporary "~" marker (this avoids ever collapsing a
#    range exactly at the story's end, which the host mishandles).
$tail = $d.Content
$tail.Find.Execute("UI~", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.MoveEnd(1, -1) | Out-Null
$tail.Collapse(0)
$d.Bookmarks.Add("_GoBack", $tail) | Out-Null

# 4) Remove the temporary marker now that the bookmark is anchored.
$d.Content.Find.Execute("~", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
